$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 8
$ws.Cells.Item(8, 1).Value = 10
$ws.Cells.Item(8, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(8, 3).Value = "La Araucanía"
$ws.Cells.Item(8, 4).Value = 44511
$ws.Cells.Item(8, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(8, 5).Value = 9
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100104
$ws.Cells.Item(8, 8).Value = "Frutos de pepita"
$ws.Cells.Item(8, 9).Value = 100104004
$ws.Cells.Item(8, 10).Value = "Níspero"
$ws.Cells.Item(8, 11).Value = "Californiana(o)"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 45
$ws.Cells.Item(8, 14).Value = 28000
$ws.Cells.Item(8, 15).Value = 28000
$ws.Cells.Item(8, 16).Value = 28000
$ws.Cells.Item(8, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(8, 19).Value = 2800
$ws.Cells.Item(8, 20).Value = 10

# New row 9
$ws.Cells.Item(9, 1).Value = 10
$ws.Cells.Item(9, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(9, 3).Value = "La Araucanía"
$ws.Cells.Item(9, 4).Value = 44511
$ws.Cells.Item(9, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100104
$ws.Cells.Item(9, 8).Value = "Frutos de pepita"
$ws.Cells.Item(9, 9).Value = 100104004
$ws.Cells.Item(9, 10).Value = "Níspero"
$ws.Cells.Item(9, 11).Value = "Californiana(o)"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 45
$ws.Cells.Item(9, 14).Value = 3200
$ws.Cells.Item(9, 15).Value = 3200
$ws.Cells.Item(9, 16).Value = 3200
$ws.Cells.Item(9, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(9, 19).Value = 320
$ws.Cells.Item(9, 20).Value = 10
